# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) / "Valor Mora" (column F) block in rows
# 16-27 is reordered from descending period order (1712 down to 1701) to
# ascending chronological order (1701 up to 1712), carrying the special
# "Valor Mora" figure (12039, vs. the usual 32834) along with period 1712
# instead of period 1712 being first.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 16
$periods  = @("1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712")
$valores  = @(32834,32834,32834,32834,32834,32834,32834,32834,32834,32834,32834,12039)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("E$row").Value = $periods[$i]
    $ws.Range("F$row").Value = $valores[$i]
}

Write-Host "Updated periods 1701-1712 across rows $startRow-$($startRow + $periods.Length - 1)"
